$wb = $excel.ActiveWorkbook

# --- Typography sheet: set G4/H4 to "0-9" ---
$wsTypo = $wb.Worksheets.Item("Typography")
$wsTypo.Range("G4").Value = "0-9"
$wsTypo.Range("H4").Value = "0-9"

# --- Translation sheet: fill rows 4-6 ---
$wsTrans = $wb.Worksheets.Item("Translation")

$wsTrans.Range("B4").Value = "SingleUseId1"
$wsTrans.Range("C4").Value = "Default"
$wsTrans.Range("D4").Value = "Left"
$wsTrans.Range("E4").Value = "LTR"
$wsTrans.Range("F4").Value = "ADC value = <value>"

$wsTrans.Range("B5").Value = "SingleUseId2"
$wsTrans.Range("C5").Value = "Default"
$wsTrans.Range("D5").Value = "Left"
$wsTrans.Range("E5").Value = "LTR"
$wsTrans.Range("F5").NumberFormat = "@"
$wsTrans.Range("F5").Value = "10"

$wsTrans.Range("B6").Value = "SingleUseId3"
$wsTrans.Range("C6").Value = "Default"
$wsTrans.Range("D6").Value = "Center"
$wsTrans.Range("E6").Value = "LTR"
$wsTrans.Range("F6").Value = "New ADC"
